$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 72038
$ws.Cells.Item(2, 4).Value = 18644
$ws.Cells.Item(2, 5).Value = 192390213
$ws.Cells.Item(3, 3).Value = 175468
$ws.Cells.Item(3, 5).Value = 582419240
$ws.Cells.Item(4, 3).Value = 68239
$ws.Cells.Item(4, 5).Value = 316086376
$ws.Cells.Item(5, 3).Value = 23991
$ws.Cells.Item(5, 5).Value = 157618645
$ws.Cells.Item(6, 3).Value = 10838
$ws.Cells.Item(6, 5).Value = 132298199
$ws.Cells.Item(7, 3).Value = 3036
$ws.Cells.Item(7, 5).Value = 88469494
$ws.Cells.Item(14, 3).Value = 76184
$ws.Cells.Item(14, 4).Value = 20386
$ws.Cells.Item(14, 5).Value = 140952636
$ws.Cells.Item(15, 3).Value = 17985
$ws.Cells.Item(15, 5).Value = 47514636
$ws.Cells.Item(16, 3).Value = 49038
$ws.Cells.Item(16, 5).Value = 160308808
$ws.Cells.Item(17, 3).Value = 17313
$ws.Cells.Item(17, 5).Value = 75498237
$ws.Cells.Item(18, 3).Value = 5449
$ws.Cells.Item(18, 5).Value = 31906552
$ws.Cells.Item(19, 3).Value = 2256
$ws.Cells.Item(19, 5).Value = 24393938
$ws.Cells.Item(23, 3).Value = 17563
$ws.Cells.Item(23, 5).Value = 31340328
$ws.Cells.Item(24, 3).Value = 25638
$ws.Cells.Item(24, 5).Value = 77436852
$ws.Cells.Item(25, 3).Value = 61149
$ws.Cells.Item(25, 5).Value = 216149543
$ws.Cells.Item(26, 3).Value = 22685
$ws.Cells.Item(26, 5).Value = 105640625
$ws.Cells.Item(27, 3).Value = 7313
$ws.Cells.Item(27, 5).Value = 45279057
$ws.Cells.Item(28, 3).Value = 2961
$ws.Cells.Item(28, 5).Value = 32872005
$ws.Cells.Item(33, 3).Value = 19226
$ws.Cells.Item(33, 5).Value = 34643529
$ws.Cells.Item(34, 3).Value = 14080
$ws.Cells.Item(34, 5).Value = 37399318
$ws.Cells.Item(35, 3).Value = 41357
$ws.Cells.Item(35, 5).Value = 130920962
$ws.Cells.Item(36, 3).Value = 15749
$ws.Cells.Item(36, 5).Value = 67212084
$ws.Cells.Item(37, 3).Value = 4642
$ws.Cells.Item(37, 4).Value = 1357
$ws.Cells.Item(37, 5).Value = 26186009
$ws.Cells.Item(38, 3).Value = 1888
$ws.Cells.Item(38, 4).Value = 569
$ws.Cells.Item(38, 5).Value = 21423365
$ws.Cells.Item(41, 3).Value = 13937
$ws.Cells.Item(41, 5).Value = 24550547
$ws.Cells.Item(42, 3).Value = 7067
$ws.Cells.Item(42, 5).Value = 23146669
$ws.Cells.Item(43, 3).Value = 16024
$ws.Cells.Item(43, 5).Value = 54989424
$ws.Cells.Item(44, 3).Value = 7152
$ws.Cells.Item(44, 5).Value = 34112846
$ws.Cells.Item(49, 3).Value = 4539
$ws.Cells.Item(49, 5).Value = 7844230
$ws.Cells.Item(50, 3).Value = 32838
$ws.Cells.Item(50, 5).Value = 88256208
$ws.Cells.Item(51, 3).Value = 98449
$ws.Cells.Item(51, 5).Value = 322268629
$ws.Cells.Item(52, 3).Value = 39809
$ws.Cells.Item(52, 4).Value = 11301
$ws.Cells.Item(52, 5).Value = 174121132
$ws.Cells.Item(53, 3).Value = 14131
$ws.Cells.Item(53, 4).Value = 3952
$ws.Cells.Item(53, 5).Value = 82495398
$ws.Cells.Item(54, 3).Value = 6058
$ws.Cells.Item(54, 5).Value = 67024902
$ws.Cells.Item(55, 3).Value = 1536
$ws.Cells.Item(55, 5).Value = 42918400
$ws.Cells.Item(59, 3).Value = 35487
$ws.Cells.Item(59, 4).Value = 11079
$ws.Cells.Item(59, 5).Value = 78863091
$ws.Cells.Item(60, 3).Value = 2980
$ws.Cells.Item(60, 5).Value = 5312026
$ws.Cells.Item(61, 3).Value = 9822
$ws.Cells.Item(61, 5).Value = 18408611
$ws.Cells.Item(64, 3).Value = 333
$ws.Cells.Item(64, 5).Value = 1199040
$ws.Cells.Item(67, 3).Value = 10728
$ws.Cells.Item(67, 5).Value = 16859614
$ws.Cells.Item(68, 3).Value = 2006
$ws.Cells.Item(68, 5).Value = 4557672
$ws.Cells.Item(72, 3).Value = 234
$ws.Cells.Item(72, 5).Value = 668722
$ws.Cells.Item(74, 3).Value = 3197
$ws.Cells.Item(74, 5).Value = 6302309
$ws.Cells.Item(75, 3).Value = 28095
$ws.Cells.Item(75, 5).Value = 68815244
$ws.Cells.Item(76, 3).Value = 83040
$ws.Cells.Item(76, 4).Value = 22423
$ws.Cells.Item(76, 5).Value = 255444506
$ws.Cells.Item(77, 3).Value = 32169
$ws.Cells.Item(77, 5).Value = 137408612
$ws.Cells.Item(78, 3).Value = 10798
$ws.Cells.Item(78, 5).Value = 58730855
$ws.Cells.Item(79, 3).Value = 4384
$ws.Cells.Item(79, 4).Value = 1281
$ws.Cells.Item(79, 5).Value = 46780990
$ws.Cells.Item(86, 3).Value = 25606
$ws.Cells.Item(86, 5).Value = 44782870
$ws.Cells.Item(87, 3).Value = 104424
$ws.Cells.Item(87, 4).Value = 24173
$ws.Cells.Item(87, 5).Value = 267667094
$ws.Cells.Item(88, 3).Value = 282734
$ws.Cells.Item(88, 4).Value = 69387
$ws.Cells.Item(88, 5).Value = 845704643
$ws.Cells.Item(89, 3).Value = 135657
$ws.Cells.Item(89, 4).Value = 33327
$ws.Cells.Item(89, 5).Value = 594250011
$ws.Cells.Item(90, 3).Value = 55029
$ws.Cells.Item(90, 4).Value = 13238
$ws.Cells.Item(90, 5).Value = 346324285
$ws.Cells.Item(91, 3).Value = 24982
$ws.Cells.Item(91, 4).Value = 6516
$ws.Cells.Item(91, 5).Value = 292086384
$ws.Cells.Item(92, 3).Value = 6727
$ws.Cells.Item(92, 5).Value = 202045011
$ws.Cells.Item(100, 3).Value = 97004
$ws.Cells.Item(100, 4).Value = 22387
$ws.Cells.Item(100, 5).Value = 173105826
$ws.Cells.Item(101, 3).Value = 6223
$ws.Cells.Item(101, 5).Value = 10719412
$ws.Cells.Item(102, 3).Value = 14575
$ws.Cells.Item(102, 5).Value = 25891821
$ws.Cells.Item(103, 3).Value = 4626
$ws.Cells.Item(103, 5).Value = 9153499
$ws.Cells.Item(104, 3).Value = 1665
$ws.Cells.Item(104, 5).Value = 4090291
$ws.Cells.Item(108, 3).Value = 6952
$ws.Cells.Item(108, 5).Value = 9819642
$ws.Cells.Item(109, 3).Value = 2583
$ws.Cells.Item(109, 4).Value = 683
$ws.Cells.Item(109, 5).Value = 5331232
$ws.Cells.Item(110, 3).Value = 8071
$ws.Cells.Item(110, 5).Value = 18197744
$ws.Cells.Item(112, 3).Value = 1093
$ws.Cells.Item(112, 5).Value = 3411096
$ws.Cells.Item(113, 3).Value = 401
$ws.Cells.Item(113, 5).Value = 1955267
$ws.Cells.Item(114, 3).Value = 99
$ws.Cells.Item(114, 4).Value = 39
$ws.Cells.Item(114, 5).Value = 1018363
$ws.Cells.Item(116, 3).Value = 5586
$ws.Cells.Item(116, 5).Value = 8454575
$ws.Cells.Item(123, 3).Value = 20435
$ws.Cells.Item(123, 5).Value = 56523499
$ws.Cells.Item(124, 3).Value = 54435
$ws.Cells.Item(124, 5).Value = 179027099
$ws.Cells.Item(125, 3).Value = 20524
$ws.Cells.Item(125, 5).Value = 90753195
$ws.Cells.Item(126, 3).Value = 6910
$ws.Cells.Item(126, 5).Value = 40877614
$ws.Cells.Item(127, 3).Value = 2719
$ws.Cells.Item(127, 5).Value = 30371856
$ws.Cells.Item(131, 3).Value = 16032
$ws.Cells.Item(131, 5).Value = 28093615
$ws.Cells.Item(132, 3).Value = 57246
$ws.Cells.Item(132, 5).Value = 160594093
$ws.Cells.Item(133, 3).Value = 120740
$ws.Cells.Item(133, 5).Value = 385082778
$ws.Cells.Item(134, 3).Value = 43605
$ws.Cells.Item(134, 5).Value = 188436643
$ws.Cells.Item(135, 3).Value = 15024
$ws.Cells.Item(135, 5).Value = 86942841
$ws.Cells.Item(136, 3).Value = 6172
$ws.Cells.Item(136, 5).Value = 69085150
$ws.Cells.Item(137, 3).Value = 1857
$ws.Cells.Item(137, 5).Value = 53845985
$ws.Cells.Item(142, 3).Value = 42293
$ws.Cells.Item(142, 5).Value = 73175281
$ws.Cells.Item(143, 3).Value = 70339
$ws.Cells.Item(143, 5).Value = 198449029
$ws.Cells.Item(144, 3).Value = 144178
$ws.Cells.Item(144, 4).Value = 38443
$ws.Cells.Item(144, 5).Value = 451535853
$ws.Cells.Item(145, 3).Value = 49222
$ws.Cells.Item(145, 5).Value = 213535064
$ws.Cells.Item(146, 3).Value = 15837
$ws.Cells.Item(146, 5).Value = 90973551
$ws.Cells.Item(147, 3).Value = 6340
$ws.Cells.Item(147, 5).Value = 67755832
$ws.Cells.Item(148, 3).Value = 1712
$ws.Cells.Item(148, 4).Value = 618
$ws.Cells.Item(148, 5).Value = 50638118
$ws.Cells.Item(153, 3).Value = 53300
$ws.Cells.Item(153, 4).Value = 15119
$ws.Cells.Item(153, 5).Value = 88404099
$ws.Cells.Item(154, 3).Value = 24722
$ws.Cells.Item(154, 5).Value = 69259679
$ws.Cells.Item(155, 3).Value = 60158
$ws.Cells.Item(155, 4).Value = 16821
$ws.Cells.Item(155, 5).Value = 202507592
$ws.Cells.Item(156, 3).Value = 24099
$ws.Cells.Item(156, 5).Value = 107796976
$ws.Cells.Item(157, 3).Value = 7505
$ws.Cells.Item(157, 5).Value = 44121629
$ws.Cells.Item(162, 3).Value = 17491
$ws.Cells.Item(162, 5).Value = 30655355
$ws.Cells.Item(163, 3).Value = 69132
$ws.Cells.Item(163, 5).Value = 192629143
$ws.Cells.Item(164, 3).Value = 160590
$ws.Cells.Item(164, 4).Value = 40206
$ws.Cells.Item(164, 5).Value = 516645773
$ws.Cells.Item(165, 3).Value = 54332
$ws.Cells.Item(165, 5).Value = 248791171
$ws.Cells.Item(166, 3).Value = 17877
$ws.Cells.Item(166, 5).Value = 115162646
$ws.Cells.Item(167, 3).Value = 7880
$ws.Cells.Item(167, 5).Value = 92124643
$ws.Cells.Item(168, 3).Value = 2217
$ws.Cells.Item(168, 5).Value = 67224808
$ws.Cells.Item(173, 3).Value = 50649
$ws.Cells.Item(173, 5).Value = 87453198
